$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Rename "Levin Lee" -> "Kevin Lee" ---
$ws.Range("F41").Value = "Kevin Lee"
$ws.Range("F52").Value = "Kevin Lee"
$ws.Range("H60").Value = "Kevin Lee"
$ws.Range("H64").Value = "Kevin Lee"
$ws.Range("F66").Value = "Kevin Lee"
$ws.Range("H67").Value = "Kevin Lee"
$ws.Range("G71").Value = "Kevin Lee"
$ws.Range("I75").Value = "Kevin Lee"
$ws.Range("I79").Value = "Kevin Lee"
$ws.Range("G80").Value = "Kevin Lee"
$ws.Range("I81").Value = "Kevin Lee"
$ws.Range("G86").Value = "Kevin Lee"
$ws.Range("G91").Value = "Kevin Lee"
$ws.Range("G95").Value = "Kevin Lee"
$ws.Range("G96").Value = "Kevin Lee"
$ws.Range("I97").Value = "Kevin Lee"
$ws.Range("I98").Value = "Kevin Lee"

# --- Step 2: Rename "Sam Carswell-Tellis" -> "Sam Tellis" ---
$ws.Range("H38").Value = "Sam Tellis"
$ws.Range("F44").Value = "Sam Tellis"

# --- Step 3: Add new rows 99-112, copy style (bold/border) from column A of existing data rows ---
$ws.Range("A98").Copy($ws.Range("A99:A112"))

# Row 99
$ws.Range("A99").Value = 97
$ws.Range("B99").Value = "2023_2s"
$ws.Range("C99").Value = 1
$ws.Range("D99").Value = $false
$ws.Range("E99").Value = $false
$ws.Range("F99").Value = "Rohan Chowla"
$ws.Range("G99").Value = "Kevin Lee"
$ws.Range("H99").Value = "Noah Dale"
$ws.Range("I99").Value = "Yafu LastName"
$ws.Range("J99").Value = 6
$ws.Range("K99").Value = 5
$ws.Range("L99").Value = 1
$ws.Range("M99").Value = 8
$ws.Range("N99").Value = 6
$ws.Range("O99").Value = -7
$ws.Range("P99").Value = 0
$ws.Range("Q99").Value = 50
$ws.Range("R99").Value = 50
$ws.Range("S99").Value = -1
$ws.Range("T99").Value = -0.1666666666666667
$ws.Range("U99").Value = -0.1666666666666667
$ws.Range("V99").Value = 0.1666666666666667

# Row 100
$ws.Range("A100").Value = 98
$ws.Range("B100").Value = "2023_2s"
$ws.Range("C100").Value = 2
$ws.Range("D100").Value = $false
$ws.Range("E100").Value = $false
$ws.Range("F100").Value = "Gabe Silverstein"
$ws.Range("G100").Value = "Alex LastName"
$ws.Range("H100").Value = "Luci Nguyen"
$ws.Range("I100").Value = "Matthew Rusten"
$ws.Range("J100").Value = 3
$ws.Range("K100").Value = 6
$ws.Range("L100").Value = 4
$ws.Range("M100").Value = 5
$ws.Range("N100").Value = 6
$ws.Range("O100").Value = -1
$ws.Range("P100").Value = 0
$ws.Range("Q100").Value = 50
$ws.Range("R100").Value = 50
$ws.Range("S100").Value = 3
$ws.Range("T100").Value = 0.5
$ws.Range("U100").Value = 0.5
$ws.Range("V100").Value = 0.5

# Row 101
$ws.Range("A101").Value = 99
$ws.Range("B101").Value = "2023_2s"
$ws.Range("C101").Value = 3
$ws.Range("D101").Value = $false
$ws.Range("E101").Value = $false
$ws.Range("F101").Value = "Cason Duszak"
$ws.Range("G101").Value = "Eric LastName"
$ws.Range("H101").Value = "Piper Parker"
$ws.Range("I101").Value = "Kim LastName"
$ws.Range("J101").Value = 6
$ws.Range("K101").Value = 3
$ws.Range("L101").Value = 2
$ws.Range("M101").Value = 7
$ws.Range("N101").Value = 6
$ws.Range("O101").Value = -5
$ws.Range("P101").Value = 0
$ws.Range("Q101").Value = 50
$ws.Range("R101").Value = 50
$ws.Range("S101").Value = -3
$ws.Range("T101").Value = -0.5
$ws.Range("U101").Value = -0.5
$ws.Range("V101").Value = 0.5

# Row 102
$ws.Range("A102").Value = 100
$ws.Range("B102").Value = "2023_2s"
$ws.Range("C102").Value = 4
$ws.Range("D102").Value = $false
$ws.Range("E102").Value = $false
$ws.Range("F102").Value = "Nathan Snow"
$ws.Range("G102").Value = "Jason Jackson"
$ws.Range("H102").Value = "Julie Jackson"
$ws.Range("I102").Value = "Carolyn LastName"
$ws.Range("J102").Value = 6
$ws.Range("K102").Value = 4
$ws.Range("L102").Value = 3
$ws.Range("M102").Value = 6
$ws.Range("N102").Value = 6
$ws.Range("O102").Value = -3
$ws.Range("P102").Value = 0
$ws.Range("Q102").Value = 50
$ws.Range("R102").Value = 50
$ws.Range("S102").Value = -2
$ws.Range("T102").Value = -0.3333333333333333
$ws.Range("U102").Value = -0.3333333333333333
$ws.Range("V102").Value = 0.3333333333333333

# Row 103
$ws.Range("A103").Value = 101
$ws.Range("B103").Value = "2023_2s"
$ws.Range("C103").Value = 5
$ws.Range("D103").Value = $true
$ws.Range("E103").Value = $true
$ws.Range("F103").Value = "Noah Dale"
$ws.Range("G103").Value = "Yafu LastName"
$ws.Range("H103").Value = "Gabe Silverstein"
$ws.Range("I103").Value = "Alex LastName"
$ws.Range("J103").Value = 5
$ws.Range("K103").Value = 6
$ws.Range("L103").Value = 8
$ws.Range("M103").Value = 4
$ws.Range("N103").Value = 6
$ws.Range("O103").Value = 4
$ws.Range("P103").Value = 0
$ws.Range("Q103").Value = 50
$ws.Range("R103").Value = 50
$ws.Range("S103").Value = 1
$ws.Range("T103").Value = 0.1666666666666667
$ws.Range("U103").Value = 0.1666666666666667
$ws.Range("V103").Value = 0.1666666666666667

# Row 104
$ws.Range("A104").Value = 102
$ws.Range("B104").Value = "2023_2s"
$ws.Range("C104").Value = 6
$ws.Range("D104").Value = $true
$ws.Range("E104").Value = $true
$ws.Range("F104").Value = "Piper Parker"
$ws.Range("G104").Value = "Kim LastName"
$ws.Range("H104").Value = "Julie Jackson"
$ws.Range("I104").Value = "Carolyn LastName"
$ws.Range("J104").Value = 2
$ws.Range("K104").Value = 6
$ws.Range("L104").Value = 7
$ws.Range("M104").Value = 6
$ws.Range("N104").Value = 6
$ws.Range("O104").Value = 1
$ws.Range("P104").Value = 0
$ws.Range("Q104").Value = 50
$ws.Range("R104").Value = 50
$ws.Range("S104").Value = 4
$ws.Range("T104").Value = 0.6666666666666666
$ws.Range("U104").Value = 0.6666666666666666
$ws.Range("V104").Value = 0.6666666666666666

# Row 105
$ws.Range("A105").Value = 103
$ws.Range("B105").Value = "2023_2s"
$ws.Range("C105").Value = 7
$ws.Range("D105").Value = $false
$ws.Range("E105").Value = $false
$ws.Range("F105").Value = "Rohan Chowla"
$ws.Range("G105").Value = "Kevin Lee"
$ws.Range("H105").Value = "Luci Nguyen"
$ws.Range("I105").Value = "Matthew Rusten"
$ws.Range("J105").Value = 6
$ws.Range("K105").Value = 2
$ws.Range("L105").Value = 1
$ws.Range("M105").Value = 5
$ws.Range("N105").Value = 6
$ws.Range("O105").Value = -4
$ws.Range("P105").Value = 0
$ws.Range("Q105").Value = 50
$ws.Range("R105").Value = 50
$ws.Range("S105").Value = -4
$ws.Range("T105").Value = -0.6666666666666666
$ws.Range("U105").Value = -0.6666666666666666
$ws.Range("V105").Value = 0.6666666666666666

# Row 106
$ws.Range("A106").Value = 104
$ws.Range("B106").Value = "2023_2s"
$ws.Range("C106").Value = 8
$ws.Range("D106").Value = $false
$ws.Range("E106").Value = $false
$ws.Range("F106").Value = "Cason Duszak"
$ws.Range("G106").Value = "Eric LastName"
$ws.Range("H106").Value = "Nathan Snow"
$ws.Range("I106").Value = "Jason Jackson"
$ws.Range("J106").Value = 6
$ws.Range("K106").Value = 5
$ws.Range("L106").Value = 2
$ws.Range("M106").Value = 3
$ws.Range("N106").Value = 6
$ws.Range("O106").Value = -1
$ws.Range("P106").Value = 0
$ws.Range("Q106").Value = 50
$ws.Range("R106").Value = 50
$ws.Range("S106").Value = -1
$ws.Range("T106").Value = -0.1666666666666667
$ws.Range("U106").Value = -0.1666666666666667
$ws.Range("V106").Value = 0.1666666666666667

# Row 107
$ws.Range("A107").Value = 105
$ws.Range("B107").Value = "2023_2s"
$ws.Range("C107").Value = 9
$ws.Range("D107").Value = $true
$ws.Range("E107").Value = $true
$ws.Range("F107").Value = "Luci Nguyen"
$ws.Range("G107").Value = "Matthew Rusten"
$ws.Range("H107").Value = "Julie Jackson"
$ws.Range("I107").Value = "Carolyn LastName"
$ws.Range("J107").Value = 6
$ws.Range("K107").Value = 5
$ws.Range("L107").Value = 5
$ws.Range("M107").Value = 6
$ws.Range("N107").Value = 6
$ws.Range("O107").Value = -1
$ws.Range("P107").Value = 0
$ws.Range("Q107").Value = 50
$ws.Range("R107").Value = 50
$ws.Range("S107").Value = -1
$ws.Range("T107").Value = -0.1666666666666667
$ws.Range("U107").Value = -0.1666666666666667
$ws.Range("V107").Value = 0.1666666666666667

# Row 108
$ws.Range("A108").Value = 106
$ws.Range("B108").Value = "2023_2s"
$ws.Range("C108").Value = 10
$ws.Range("D108").Value = $true
$ws.Range("E108").Value = $true
$ws.Range("F108").Value = "Nathan Snow"
$ws.Range("G108").Value = "Jason Jackson"
$ws.Range("H108").Value = "Gabe Silverstein"
$ws.Range("I108").Value = "Alex LastName"
$ws.Range("J108").Value = 6
$ws.Range("K108").Value = 5
$ws.Range("L108").Value = 3
$ws.Range("M108").Value = 4
$ws.Range("N108").Value = 6
$ws.Range("O108").Value = -1
$ws.Range("P108").Value = 0
$ws.Range("Q108").Value = 50
$ws.Range("R108").Value = 50
$ws.Range("S108").Value = -1
$ws.Range("T108").Value = -0.1666666666666667
$ws.Range("U108").Value = -0.1666666666666667
$ws.Range("V108").Value = 0.1666666666666667

# Row 109
$ws.Range("A109").Value = 107
$ws.Range("B109").Value = "2023_2s"
$ws.Range("C109").Value = 11
$ws.Range("D109").Value = $true
$ws.Range("E109").Value = $true
$ws.Range("F109").Value = "Nathan Snow"
$ws.Range("G109").Value = "Jason Jackson"
$ws.Range("H109").Value = "Luci Nguyen"
$ws.Range("I109").Value = "Matthew Rusten"
$ws.Range("J109").Value = 9
$ws.Range("K109").Value = 7
$ws.Range("L109").Value = 3
$ws.Range("M109").Value = 5
$ws.Range("N109").Value = 9
$ws.Range("O109").Value = -2
$ws.Range("P109").Value = 0
$ws.Range("Q109").Value = 50
$ws.Range("R109").Value = 50
$ws.Range("S109").Value = -2
$ws.Range("T109").Value = -0.3333333333333333
$ws.Range("U109").Value = -0.5
$ws.Range("V109").Value = 0.5

# Row 110
$ws.Range("A110").Value = 108
$ws.Range("B110").Value = "2023_2s"
$ws.Range("C110").Value = 12
$ws.Range("D110").Value = $false
$ws.Range("E110").Value = $false
$ws.Range("F110").Value = "Rohan Chowla"
$ws.Range("G110").Value = "Kevin Lee"
$ws.Range("H110").Value = "Cason Duszak"
$ws.Range("I110").Value = "Eric LastName"
$ws.Range("J110").Value = 6
$ws.Range("K110").Value = 5
$ws.Range("L110").Value = 1
$ws.Range("M110").Value = 2
$ws.Range("N110").Value = 6
$ws.Range("O110").Value = -1
$ws.Range("P110").Value = 0
$ws.Range("Q110").Value = 50
$ws.Range("R110").Value = 50
$ws.Range("S110").Value = -1
$ws.Range("T110").Value = -0.1666666666666667
$ws.Range("U110").Value = -0.1666666666666667
$ws.Range("V110").Value = 0.1666666666666667

# Row 111
$ws.Range("A111").Value = 109
$ws.Range("B111").Value = "2023_2s"
$ws.Range("C111").Value = 13
$ws.Range("D111").Value = $true
$ws.Range("E111").Value = $true
$ws.Range("F111").Value = "Cason Duszak"
$ws.Range("G111").Value = "Eric LastName"
$ws.Range("H111").Value = "Nathan Snow"
$ws.Range("I111").Value = "Jason Jackson"
$ws.Range("J111").Value = 6
$ws.Range("K111").Value = 5
$ws.Range("L111").Value = 2
$ws.Range("M111").Value = 3
$ws.Range("N111").Value = 6
$ws.Range("O111").Value = -1
$ws.Range("P111").Value = 0
$ws.Range("Q111").Value = 50
$ws.Range("R111").Value = 50
$ws.Range("S111").Value = -1
$ws.Range("T111").Value = -0.1666666666666667
$ws.Range("U111").Value = -0.1666666666666667
$ws.Range("V111").Value = 0.1666666666666667

# Row 112
$ws.Range("A112").Value = 110
$ws.Range("B112").Value = "2023_2s"
$ws.Range("C112").Value = 14
$ws.Range("D112").Value = $false
$ws.Range("E112").Value = $true
$ws.Range("F112").Value = "Rohan Chowla"
$ws.Range("G112").Value = "Kevin Lee"
$ws.Range("H112").Value = "Cason Duszak"
$ws.Range("I112").Value = "Eric LastName"
$ws.Range("J112").Value = 9
$ws.Range("K112").Value = 8
$ws.Range("L112").Value = 1
$ws.Range("M112").Value = 2
$ws.Range("N112").Value = 9
$ws.Range("O112").Value = -1
$ws.Range("P112").Value = 0
$ws.Range("Q112").Value = 50
$ws.Range("R112").Value = 50
$ws.Range("S112").Value = -1
$ws.Range("T112").Value = -0.1666666666666667
$ws.Range("U112").Value = -0.25
$ws.Range("V112").Value = 0.25

